$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1802.8658
$ws.Range("I15").Value = 1802.8658
$ws.Range("K15").Value = 5408.597400000001
$ws.Range("M15").Value = -5239.597400000001
$ws.Range("H33").Value = 58925.65
$ws.Range("I33").Value = 66769.07000000001
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 66769.07000000001
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -66540.07000000001
$ws.Range("N33").Value = -558
$ws.Range("H64").Value = 27500.75
$ws.Range("I64").Value = 51650
$ws.Range("J64").Value = 3351.5
$ws.Range("K64").Value = 51650
$ws.Range("L64").Value = 3351.5
$ws.Range("M64").Value = -51402
$ws.Range("N64").Value = -3847.5
$ws.Range("H67").Value = 27500.75
$ws.Range("I67").Value = 51650
$ws.Range("J67").Value = 3351.5
$ws.Range("K67").Value = 51650
$ws.Range("L67").Value = 3351.5
$ws.Range("M67").Value = -50792
$ws.Range("N67").Value = -5067.5
$ws.Range("H132").Value = 2622.9119
$ws.Range("I132").Value = 2224.125
$ws.Range("J132").Value = 3580
$ws.Range("K132").Value = 6672.375
$ws.Range("L132").Value = 10740
$ws.Range("M132").Value = -4142.375
$ws.Range("N132").Value = -15800
$ws.Range("H137").Value = 11365663
$ws.Range("I137").Value = 1919.75
$ws.Range("J137").Value = 20835448
$ws.Range("K137").Value = 5759.25
$ws.Range("L137").Value = 62506344
$ws.Range("M137").Value = -3209.25
$ws.Range("N137").Value = -62511444

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21755.578
$ws.Range("I32").Value = 20303.117
$ws.Range("J32").Value = 34101.5
$ws.Range("K32").Value = 20303.117
$ws.Range("L32").Value = 34101.5
$ws.Range("M32").Value = -20016.117
$ws.Range("N32").Value = -34675.5
$ws.Range("H74").Value = 845.5217
$ws.Range("I74").Value = 900.2273
$ws.Range("J74").Value = 795.375
$ws.Range("K74").Value = 900.2273
$ws.Range("L74").Value = 795.375
$ws.Range("M74").Value = -26.22730000000001
$ws.Range("N74").Value = -2543.375
$ws.Range("H77").Value = 845.5217
$ws.Range("I77").Value = 900.2273
$ws.Range("J77").Value = 795.375
$ws.Range("K77").Value = 4501.136500000001
$ws.Range("L77").Value = 3976.875
$ws.Range("M77").Value = -133.1365000000005
$ws.Range("N77").Value = -12712.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 60000
$ws.Range("J50").Value = 60000
$ws.Range("L50").Value = 60000
$ws.Range("N50").Value = -61148

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1765.9231
$ws.Range("I31").Value = 974.1177
$ws.Range("J31").Value = 2150.5144
$ws.Range("K31").Value = 974.1177
$ws.Range("L31").Value = 2150.5144
$ws.Range("M31").Value = -679.1177
$ws.Range("N31").Value = -2740.5144
$ws.Range("H34").Value = 1765.9231
$ws.Range("I34").Value = 974.1177
$ws.Range("J34").Value = 2150.5144
$ws.Range("K34").Value = 974.1177
$ws.Range("L34").Value = 2150.5144
$ws.Range("M34").Value = -772.1177
$ws.Range("N34").Value = -2554.5144
$ws.Range("H58").Value = 4852.8
$ws.Range("I58").Value = 1194.1818
$ws.Range("J58").Value = 6970.9473
$ws.Range("K58").Value = 1194.1818
$ws.Range("L58").Value = 6970.9473
$ws.Range("M58").Value = -991.1818000000001
$ws.Range("N58").Value = -7376.9473
$ws.Range("H62").Value = 11092.083
$ws.Range("I62").Value = 11943.571
$ws.Range("J62").Value = 9900
$ws.Range("K62").Value = 11943.571
$ws.Range("L62").Value = 9900
$ws.Range("M62").Value = -11319.571
$ws.Range("N62").Value = -11148
$ws.Range("H65").Value = 11092.083
$ws.Range("I65").Value = 11943.571
$ws.Range("J65").Value = 9900
$ws.Range("K65").Value = 59717.855
$ws.Range("L65").Value = 49500
$ws.Range("M65").Value = -56597.855
$ws.Range("N65").Value = -55740
$ws.Range("H136").Value = 4852.8
$ws.Range("I136").Value = 1194.1818
$ws.Range("J136").Value = 6970.9473
$ws.Range("K136").Value = 3582.5454
$ws.Range("L136").Value = 20912.8419
$ws.Range("M136").Value = -1032.5454
$ws.Range("N136").Value = -26012.8419
$ws.Range("H140").Value = 46120
$ws.Range("J140").Value = 46120
$ws.Range("L140").Value = 46120
$ws.Range("N140").Value = -56480

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 5730
$ws.Range("J51").Value = 5730
$ws.Range("L51").Value = 17190
$ws.Range("N51").Value = -18110
$ws.Range("H58").Value = 2001.6666
$ws.Range("I58").Value = 2001.6666
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 6004.9998
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5876.9998
$ws.Range("H68").Value = 783.2447
$ws.Range("I68").Value = 526.59015
$ws.Range("J68").Value = 1257.6666
$ws.Range("K68").Value = 1579.77045
$ws.Range("L68").Value = 3772.9998
$ws.Range("M68").Value = -768.77045
$ws.Range("N68").Value = -5394.9998
$ws.Range("H71").Value = 783.2447
$ws.Range("I71").Value = 526.59015
$ws.Range("J71").Value = 1257.6666
$ws.Range("K71").Value = 4739.31135
$ws.Range("L71").Value = 11318.9994
$ws.Range("M71").Value = -683.3113499999999
$ws.Range("N71").Value = -19430.9994
$ws.Range("H97").Value = 933.3333
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 3600
$ws.Range("M97").Value = -704
$ws.Range("N97").Value = -4592
$ws.Range("H107").Value = 195719.7
$ws.Range("I107").Value = 451.6111
$ws.Range("J107").Value = 355484.5
$ws.Range("K107").Value = 1354.8333
$ws.Range("L107").Value = 1066453.5
$ws.Range("M107").Value = 565.1667
$ws.Range("N107").Value = -1070293.5
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1389.8334
$ws.Range("I122").Value = 1332.2307
$ws.Range("J122").Value = 1539.6
$ws.Range("K122").Value = 3996.6921
$ws.Range("L122").Value = 4618.799999999999
$ws.Range("M122").Value = -1546.6921
$ws.Range("N122").Value = -9518.799999999999
$ws.Range("H126").Value = 3439.4285
$ws.Range("I126").Value = 3778
$ws.Range("K126").Value = 11334
$ws.Range("M126").Value = -8864
$ws.Range("H132").Value = 2518.5833
$ws.Range("I132").Value = 1795.7646
$ws.Range("J132").Value = 4274
$ws.Range("K132").Value = 5387.293799999999
$ws.Range("L132").Value = 12822
$ws.Range("M132").Value = -2857.293799999999
$ws.Range("N132").Value = -17882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4582.875
$ws.Range("I122").Value = 5108.28
$ws.Range("K122").Value = 15324.84
$ws.Range("M122").Value = -12874.84
$ws.Range("H132").Value = 2340.2666
$ws.Range("I132").Value = 2321.4707
$ws.Range("J132").Value = 2398.3635
$ws.Range("K132").Value = 6964.4121
$ws.Range("L132").Value = 7195.0905
$ws.Range("M132").Value = -4434.4121
$ws.Range("N132").Value = -12255.0905
$ws.Range("H136").Value = 2266.3333
$ws.Range("I136").Value = 1515.5807
$ws.Range("J136").Value = 3278.2173
$ws.Range("K136").Value = 4546.742099999999
$ws.Range("L136").Value = 9834.651899999999
$ws.Range("M136").Value = -1996.742099999999
$ws.Range("N136").Value = -14934.6519

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3522.5
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 3691.2222
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 11073.6666
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -16133.6666
$ws.Range("H136").Value = 2069.4753
$ws.Range("I136").Value = 2201.3713
$ws.Range("K136").Value = 6604.113899999999
$ws.Range("M136").Value = -4054.113899999999
